$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7389.6
$ws.Range("I43").Value = 9000
$ws.Range("J43").Value = 4974
$ws.Range("K43").Value = 9000
$ws.Range("L43").Value = 4974
$ws.Range("M43").Value = -8931
$ws.Range("N43").Value = -5112
$ws.Range("H92").Value = 639.6
$ws.Range("I92").Value = 671.5789
$ws.Range("K92").Value = 671.5789
$ws.Range("M92").Value = 576.4211
$ws.Range("H132").Value = 3971
$ws.Range("I132").Value = 3971
$ws.Range("K132").Value = 11913
$ws.Range("M132").Value = -9383
$ws.Range("H137").Value = 2549.75
$ws.Range("I137").Value = 1914
$ws.Range("K137").Value = 5742
$ws.Range("M137").Value = -3192

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6080.1
$ws.Range("I45").Value = 7655.857
$ws.Range("K45").Value = 7655.857
$ws.Range("M45").Value = -7278.857
$ws.Range("H76").Value = 52499
$ws.Range("J76").Value = 52499
$ws.Range("L76").Value = 52499
$ws.Range("N76").Value = -53175
$ws.Range("H79").Value = 52499
$ws.Range("J79").Value = 52499
$ws.Range("L79").Value = 52499
$ws.Range("N79").Value = -54839
$ws.Range("H88").Value = 49582.57
$ws.Range("I88").Value = 112721.89
$ws.Range("J88").Value = 2228.0833
$ws.Range("K88").Value = 112721.89
$ws.Range("L88").Value = 2228.0833
$ws.Range("M88").Value = -112315.89
$ws.Range("N88").Value = -3040.0833
$ws.Range("H91").Value = 49582.57
$ws.Range("I91").Value = 112721.89
$ws.Range("J91").Value = 2228.0833
$ws.Range("K91").Value = 112721.89
$ws.Range("L91").Value = 2228.0833
$ws.Range("M91").Value = -111317.89
$ws.Range("N91").Value = -5036.0833
$ws.Range("H95").Value = 100000
$ws.Range("J95").Value = 100000
$ws.Range("L95").Value = 100000
$ws.Range("N95").Value = -105492
$ws.Range("H97").Value = 518.5
$ws.Range("I97").Value = 497.5
$ws.Range("K97").Value = 497.5
$ws.Range("M97").Value = -1.5
$ws.Range("H122").Value = 4457.9585
$ws.Range("I122").Value = 3556.5557
$ws.Range("J122").Value = 7162.1665
$ws.Range("K122").Value = 10669.6671
$ws.Range("L122").Value = 21486.4995
$ws.Range("M122").Value = -8219.667099999999
$ws.Range("N122").Value = -26386.4995

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2095.5715
$ws.Range("I86").Value = 2021.5385
$ws.Range("K86").Value = 2021.5385
$ws.Range("M86").Value = -898.5385000000001
$ws.Range("H89").Value = 2095.5715
$ws.Range("I89").Value = 2021.5385
$ws.Range("K89").Value = 10107.6925
$ws.Range("M89").Value = -4491.692500000001
$ws.Range("H105").Value = 2114.1365
$ws.Range("I105").Value = 1608
$ws.Range("J105").Value = 2999.875
$ws.Range("K105").Value = 1608
$ws.Range("L105").Value = 2999.875
$ws.Range("M105").Value = 139
$ws.Range("N105").Value = -6493.875

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10598
$ws.Range("I31").Value = 7688.727
$ws.Range("J31").Value = 14153.777
$ws.Range("K31").Value = 7688.727
$ws.Range("L31").Value = 14153.777
$ws.Range("M31").Value = -7393.727
$ws.Range("N31").Value = -14743.777
$ws.Range("H34").Value = 10598
$ws.Range("I34").Value = 7688.727
$ws.Range("J34").Value = 14153.777
$ws.Range("K34").Value = 7688.727
$ws.Range("L34").Value = 14153.777
$ws.Range("M34").Value = -7486.727
$ws.Range("N34").Value = -14557.777
$ws.Range("J58").Value = 2415.3635
$ws.Range("L58").Value = 2415.3635
$ws.Range("N58").Value = -2821.3635
$ws.Range("H92").Value = 100000
$ws.Range("J92").Value = 100000
$ws.Range("L92").Value = 100000
$ws.Range("N92").Value = -104992
$ws.Range("H122").Value = 3858.9375
$ws.Range("I122").Value = 3242.2
$ws.Range("K122").Value = 9726.599999999999
$ws.Range("M122").Value = -7276.599999999999
$ws.Range("J136").Value = 2415.3635
$ws.Range("L136").Value = 7246.0905
$ws.Range("N136").Value = -12346.0905

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 336269.34
$ws.Range("I5").Value = 999999
$ws.Range("J5").Value = 4404.5
$ws.Range("K5").Value = 2999997
$ws.Range("L5").Value = 13213.5
$ws.Range("M5").Value = -2999885
$ws.Range("N5").Value = -13437.5
$ws.Range("H32").Value = 33556244
$ws.Range("I32").Value = 50167332
$ws.Range("J32").Value = 334066
$ws.Range("K32").Value = 150501996
$ws.Range("L32").Value = 1002198
$ws.Range("M32").Value = -150501713
$ws.Range("N32").Value = -1002764
$ws.Range("H63").Value = 2996.6667
$ws.Range("I63").Value = 2996.6667
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 8990.000100000001
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -8241.000100000001
$ws.Range("H66").Value = 2996.6667
$ws.Range("I66").Value = 2996.6667
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 26970.0003
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -23226.0003
$ws.Range("H108").Value = 5825.2856
$ws.Range("I108").Value = 5825.2856
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 17475.8568
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -14595.8568
$ws.Range("H135").Value = 336269.34
$ws.Range("I135").Value = 999999
$ws.Range("J135").Value = 4404.5
$ws.Range("K135").Value = 8999991
$ws.Range("L135").Value = 39640.5
$ws.Range("M135").Value = -8997456
$ws.Range("N135").Value = -44710.5
$ws.Range("H138").Value = 4073.4
$ws.Range("I138").Value = 1791.6
$ws.Range("J138").Value = 6355.2
$ws.Range("K138").Value = 5374.799999999999
$ws.Range("L138").Value = 19065.6
$ws.Range("M138").Value = -234.7999999999993

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 46500.5
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 46500.5
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("M39").Value = 46500.5
$ws.Range("N39").Value = -47564.5
$ws.Range("H70").Value = 13571.429
$ws.Range("I70").Value = 9455.666999999999
$ws.Range("J70").Value = 20979.8
$ws.Range("K70").Value = 9455.666999999999
$ws.Range("L70").Value = 20979.8
$ws.Range("M70").Value = -9185.666999999999
$ws.Range("N70").Value = -21519.8
$ws.Range("H73").Value = 13571.429
$ws.Range("I73").Value = 9455.666999999999
$ws.Range("J73").Value = 20979.8
$ws.Range("K73").Value = 9455.666999999999
$ws.Range("L73").Value = 20979.8
$ws.Range("M73").Value = -8519.666999999999
$ws.Range("H97").Value = 1484.2916
$ws.Range("I97").Value = 1222.9412
$ws.Range("J97").Value = 2119
$ws.Range("K97").Value = 1222.9412
$ws.Range("L97").Value = 2119
$ws.Range("M97").Value = -726.9412
$ws.Range("N97").Value = -3111

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4149.8335
$ws.Range("I7").Value = 4178.8
$ws.Range("J7").Value = 4005
$ws.Range("K7").Value = 4178.8
$ws.Range("L7").Value = 4005
$ws.Range("M7").Value = -4066.8
$ws.Range("N7").Value = -4229
$ws.Range("H26").Value = 1450
$ws.Range("I26").Value = 1450
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1450
$ws.Range("L26").ClearContents()
$ws.Range("N26").Value = 0
$ws.Range("M26").Value = -1155
$ws.Range("H40").Value = 3495.1667
$ws.Range("I40").Value = 2694.2
$ws.Range("K40").Value = 2694.2
$ws.Range("M40").Value = -2558.2
$ws.Range("H57").Value = 26992.5
$ws.Range("H61").Value = 5346.6
$ws.Range("I61").Value = 5346.6
$ws.Range("K61").Value = 5346.6
$ws.Range("M61").Value = -5144.6
$ws.Range("H93").Value = 1362
$ws.Range("I93").Value = 775.6
$ws.Range("K93").Value = 775.6
$ws.Range("M93").Value = 472.4
$ws.Range("H101").Value = 27999.666
$ws.Range("J101").Value = 27999.666
$ws.Range("L101").Value = 27999.666
$ws.Range("N101").Value = -34489.666
$ws.Range("H113").Value = 5346.6
$ws.Range("I113").Value = 5346.6
$ws.Range("K113").Value = 5346.6
$ws.Range("M113").Value = -3176.6
$ws.Range("H126").Value = 4149.8335
$ws.Range("I126").Value = 4178.8
$ws.Range("J126").Value = 4005
$ws.Range("K126").Value = 12536.4
$ws.Range("L126").Value = 12015
$ws.Range("M126").Value = -10066.4
$ws.Range("N126").Value = -16955

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 59000
$ws.Range("J76").Value = 59000
$ws.Range("L76").Value = 59000
$ws.Range("N76").Value = -59630
$ws.Range("H79").Value = 59000
$ws.Range("J79").Value = 59000
$ws.Range("L79").Value = 59000
$ws.Range("N79").Value = -61184
$ws.Range("H122").Value = 2912.1333
$ws.Range("J122").Value = 3243.111
$ws.Range("L122").Value = 9729.332999999999
$ws.Range("N122").Value = -14629.333
$ws.Range("H126").Value = 2294.0625
$ws.Range("I126").Value = 2776.25
$ws.Range("K126").Value = 8328.75
$ws.Range("M126").Value = -5858.75
$ws.Range("H136").Value = 45456450
$ws.Range("I136").Value = 45456450
$ws.Range("K136").Value = 136369350
$ws.Range("M136").Value = -136366800
